$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 2258.9333
$ws.Cells.Item(40, 9).Value = 2093.4736
$ws.Cells.Item(40, 10).Value = 2544.7273
$ws.Cells.Item(40, 11).Value = 2093.4736
$ws.Cells.Item(40, 12).Value = 2544.7273
$ws.Cells.Item(40, 13).Value = -1918.4736
$ws.Cells.Item(40, 14).Value = -2894.7273
# Row 88
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()
# Row 91
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()
# Row 111
$ws.Cells.Item(111, 8).Value = 1600
$ws.Cells.Item(111, 9).Value = 1600
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 4800
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = -1733
$ws.Cells.Item(111, 14).ClearContents()
# Row 129
$ws.Cells.Item(129, 8).Value = 2523.5
$ws.Cells.Item(129, 9).Value = 886.8
$ws.Cells.Item(129, 10).Value = 3004.8823
$ws.Cells.Item(129, 11).Value = 2660.4
$ws.Cells.Item(129, 12).Value = 9014.6469
$ws.Cells.Item(129, 13).Value = 2339.6
$ws.Cells.Item(129, 14).Value = -19014.6469
# Row 135
$ws.Cells.Item(135, 8).Value = 523.3333
$ws.Cells.Item(135, 9).Value = 523.3333
$ws.Cells.Item(135, 11).Value = 4709.9997
$ws.Cells.Item(135, 13).Value = -2174.9997
# Row 137
$ws.Cells.Item(137, 8).Value = 1789.7
$ws.Cells.Item(137, 9).Value = 982.8333
$ws.Cells.Item(137, 11).Value = 2948.4999
$ws.Cells.Item(137, 13).Value = -398.4998999999998
# Row 138
$ws.Cells.Item(138, 8).Value = 4816.107
$ws.Cells.Item(138, 10).Value = 5393.952
$ws.Cells.Item(138, 12).Value = 16181.856
$ws.Cells.Item(138, 14).Value = -26461.856

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 620.75
$ws.Cells.Item(2, 9).Value = 494.33334
$ws.Cells.Item(2, 11).Value = 494.33334
$ws.Cells.Item(2, 13).Value = -381.33334
# Row 32
$ws.Cells.Item(32, 8).Value = 2202560.2
$ws.Cells.Item(32, 9).Value = 2695045
$ws.Cells.Item(32, 11).Value = 2695045
$ws.Cells.Item(32, 13).Value = -2694758
# Row 45
$ws.Cells.Item(45, 8).Value = 143978.58
$ws.Cells.Item(45, 9).Value = 167775
$ws.Cells.Item(45, 10).Value = 1200
$ws.Cells.Item(45, 11).Value = 167775
$ws.Cells.Item(45, 12).Value = 1200
$ws.Cells.Item(45, 13).Value = -167398
$ws.Cells.Item(45, 14).Value = -1954
# Row 74
$ws.Cells.Item(74, 8).Value = 616.6667
$ws.Cells.Item(74, 9).Value = 750
$ws.Cells.Item(74, 11).Value = 750
$ws.Cells.Item(74, 13).Value = 124
# Row 77
$ws.Cells.Item(77, 8).Value = 616.6667
$ws.Cells.Item(77, 9).Value = 750
$ws.Cells.Item(77, 11).Value = 3750
$ws.Cells.Item(77, 13).Value = 618
# Row 110
$ws.Cells.Item(110, 8).Value = 1545.1428
$ws.Cells.Item(110, 9).Value = 1545.1428
$ws.Cells.Item(110, 11).Value = 1545.1428
$ws.Cells.Item(110, 13).Value = 499.8571999999999
# Row 116
$ws.Cells.Item(116, 8).Value = 620.75
$ws.Cells.Item(116, 9).Value = 494.33334
$ws.Cells.Item(116, 11).Value = 494.33334
$ws.Cells.Item(116, 13).Value = 1799.66666
# Row 135
$ws.Cells.Item(135, 8).Value = 355000
$ws.Cells.Item(135, 10).Value = 355000
$ws.Cells.Item(135, 12).Value = 355000
$ws.Cells.Item(135, 14).Value = -365140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 620.75
$ws.Cells.Item(3, 9).Value = 494.33334
$ws.Cells.Item(3, 11).Value = 494.33334
$ws.Cells.Item(3, 13).Value = -380.33334
# Row 99
$ws.Cells.Item(99, 8).Value = 1987
$ws.Cells.Item(99, 9).Value = 1987
$ws.Cells.Item(99, 11).Value = 1987
$ws.Cells.Item(99, 13).Value = -489
# Row 107
$ws.Cells.Item(107, 8).Value = 2807.818
$ws.Cells.Item(107, 9).Value = 2866.7
$ws.Cells.Item(107, 10).Value = 2219
$ws.Cells.Item(107, 11).Value = 2866.7
$ws.Cells.Item(107, 12).Value = 2219
$ws.Cells.Item(107, 13).Value = -946.6999999999998
$ws.Cells.Item(107, 14).Value = -6059
# Row 134
$ws.Cells.Item(134, 8).Value = 2300.25
$ws.Cells.Item(134, 9).Value = 2339.3333
$ws.Cells.Item(134, 10).Value = 1714
$ws.Cells.Item(134, 11).Value = 7017.999899999999
$ws.Cells.Item(134, 12).Value = 5142
$ws.Cells.Item(134, 13).Value = -4482.999899999999
$ws.Cells.Item(134, 14).Value = -10212
# Row 138
$ws.Cells.Item(138, 8).Value = 100000
$ws.Cells.Item(138, 10).Value = 100000
$ws.Cells.Item(138, 12).Value = 100000
$ws.Cells.Item(138, 14).Value = -110280

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 102
$ws.Cells.Item(102, 8).Value = 19495
$ws.Cells.Item(102, 10).Value = 19495
$ws.Cells.Item(102, 12).Value = 19495
$ws.Cells.Item(102, 14).Value = -24363
# Row 107
$ws.Cells.Item(107, 8).Value = 1088.4166
$ws.Cells.Item(107, 9).Value = 1161.5714
$ws.Cells.Item(107, 10).Value = 986
$ws.Cells.Item(107, 11).Value = 1161.5714
$ws.Cells.Item(107, 12).Value = 986
$ws.Cells.Item(107, 13).Value = 758.4286
$ws.Cells.Item(107, 14).Value = -4826
# Row 109
$ws.Cells.Item(109, 8).Value = 56450
$ws.Cells.Item(109, 10).Value = 56450
$ws.Cells.Item(109, 12).Value = 56450
$ws.Cells.Item(109, 14).Value = -58530
# Row 132
$ws.Cells.Item(132, 8).Value = 2769.6
$ws.Cells.Item(132, 9).Value = 3007.4
$ws.Cells.Item(132, 10).Value = 2056.2
$ws.Cells.Item(132, 11).Value = 9022.200000000001
$ws.Cells.Item(132, 12).Value = 6168.599999999999
$ws.Cells.Item(132, 13).Value = -6492.200000000001
$ws.Cells.Item(132, 14).Value = -11228.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Cells.Item(39, 8).Value = 8500
$ws.Cells.Item(39, 10).Value = 8500
$ws.Cells.Item(39, 12).Value = 25500
$ws.Cells.Item(39, 14).Value = -26088

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 51586.5
$ws.Cells.Item(15, 10).Value = 51586.5
$ws.Cells.Item(15, 12).Value = 51586.5
$ws.Cells.Item(15, 14).Value = -52162.5
# Row 81
$ws.Cells.Item(81, 8).Value = 51586.5
$ws.Cells.Item(81, 10).Value = 51586.5
$ws.Cells.Item(81, 12).Value = 51586.5
$ws.Cells.Item(81, 14).Value = -53582.5
# Row 84
$ws.Cells.Item(84, 8).Value = 51586.5
$ws.Cells.Item(84, 10).Value = 51586.5
$ws.Cells.Item(84, 12).Value = 154759.5
$ws.Cells.Item(84, 14).Value = -164743.5
# Row 113
$ws.Cells.Item(113, 8).Value = 863.6667
$ws.Cells.Item(113, 9).Value = 863.6667
$ws.Cells.Item(113, 11).Value = 863.6667
$ws.Cells.Item(113, 13).Value = 1306.3333
# Row 122
$ws.Cells.Item(122, 8).Value = 2628.5715
$ws.Cells.Item(122, 9).Value = 2303.5
$ws.Cells.Item(122, 10).Value = 2758.6
$ws.Cells.Item(122, 11).Value = 6910.5
$ws.Cells.Item(122, 12).Value = 8275.799999999999
$ws.Cells.Item(122, 13).Value = -4460.5
$ws.Cells.Item(122, 14).Value = -13175.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Cells.Item(13, 8).Value = 3000
$ws.Cells.Item(13, 9).Value = 3000
$ws.Cells.Item(13, 11).Value = 3000
$ws.Cells.Item(13, 13).Value = -2860
# Row 16
$ws.Cells.Item(16, 8).Value = 840.6
$ws.Cells.Item(16, 9).Value = 856.5
$ws.Cells.Item(16, 11).Value = 856.5
$ws.Cells.Item(16, 13).Value = -686.5
# Row 92
$ws.Cells.Item(92, 8).Value = 55000
$ws.Cells.Item(92, 10).Value = 55000
$ws.Cells.Item(92, 12).Value = 55000
$ws.Cells.Item(92, 14).Value = -59992
# Row 100
$ws.Cells.Item(100, 8).Value = 3460.25
$ws.Cells.Item(100, 9).Value = 2748.125
$ws.Cells.Item(100, 10).Value = 4884.5
$ws.Cells.Item(100, 11).Value = 2748.125
$ws.Cells.Item(100, 12).Value = 4884.5
$ws.Cells.Item(100, 13).Value = -2207.125
$ws.Cells.Item(100, 14).Value = -5966.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Cells.Item(14, 8).Value = 31245.5
$ws.Cells.Item(14, 10).Value = 40000
$ws.Cells.Item(14, 12).Value = 40000
$ws.Cells.Item(14, 14).Value = -40336
# Row 24
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).ClearContents()
# Row 41
$ws.Cells.Item(41, 8).Value = 19981.4
$ws.Cells.Item(41, 9).Value = 19978
$ws.Cells.Item(41, 10).Value = 19982.25
$ws.Cells.Item(41, 11).Value = 19978
$ws.Cells.Item(41, 12).Value = 19982.25
$ws.Cells.Item(41, 13).Value = -19588
$ws.Cells.Item(41, 14).Value = -20762.25
# Row 107
$ws.Cells.Item(107, 8).Value = 639.75
$ws.Cells.Item(107, 9).Value = 498.125
$ws.Cells.Item(107, 10).Value = 923
$ws.Cells.Item(107, 11).Value = 1494.375
$ws.Cells.Item(107, 12).Value = 2769
$ws.Cells.Item(107, 13).Value = 425.625
$ws.Cells.Item(107, 14).Value = -6609
# Row 116
$ws.Cells.Item(116, 8).Value = 30500
$ws.Cells.Item(116, 10).Value = 30500
$ws.Cells.Item(116, 12).Value = 30500
$ws.Cells.Item(116, 14).Value = -39678
# Row 136
$ws.Cells.Item(136, 8).Value = 2230.1
$ws.Cells.Item(136, 9).Value = 2094.2
$ws.Cells.Item(136, 11).Value = 6282.599999999999
$ws.Cells.Item(136, 13).Value = -3732.599999999999
